$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WSM")

# Insert a new column before column D (shifts D:K -> E:L), for the new quarter's data
$ws.Columns("D").Insert()

# Copy number formats/styles from column E (the old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest quarter's figures
$ws.Cells.Item(7, 4).Value = 43401
$ws.Cells.Item(8, 4).Value = 1357000
$ws.Cells.Item(9, 4).Value = 862000
$ws.Cells.Item(10, 4).Value = 495000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 6000
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 1262600
$ws.Cells.Item(18, 4).Value = 94400
$ws.Cells.Item(20, 4).Value = -2300
$ws.Cells.Item(21, 4).Value = 139500
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 4).Value = 92100
$ws.Cells.Item(24, 4).Value = 25200
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 66900
$ws.Cells.Item(27, 4).Value = 66900
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 14600
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 2300
$ws.Cells.Item(33, 4).Value = 81500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 81500
$ws.Cells.Item(38, 4).Value = 43401
$ws.Cells.Item(41, 4).Value = 164400
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 113600
$ws.Cells.Item(44, 4).Value = 1197600
$ws.Cells.Item(45, 4).Value = 115900
$ws.Cells.Item(46, 4).Value = 1591400
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 931400
$ws.Cells.Item(49, 4).Value = 85600
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 110300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 2718800
$ws.Cells.Item(57, 4).Value = 487700
$ws.Cells.Item(58, 4).Value = 60000
$ws.Cells.Item(59, 4).Value = 489000
$ws.Cells.Item(60, 4).Value = 1036700
$ws.Cells.Item(61, 4).Value = 299600
$ws.Cells.Item(62, 4).Value = 290500
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 1626800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 532200
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 1091900
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43401
$ws.Cells.Item(81, 4).Value = 81500
$ws.Cells.Item(83, 4).Value = 47400
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 59400
$ws.Cells.Item(91, 4).Value = -48300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -47000
$ws.Cells.Item(96, 4).Value = -35300
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -22300
$ws.Cells.Item(101, 4).Value = -200
$ws.Cells.Item(102, 4).Value = -10200
